# "Generate Report for Handback"
#
# The handback sync completed, so the localization-status report is
# regenerated: the "Ready for handoff" status becomes
# "Handed back: in sync with en-US" everywhere it appears (Overview!E/F
# and the per-locale Status column), the de-de row's handoff/handback
# file + datetime now point at the real de-de xliff (mirroring zh-cn,
# which was already correct), the stale "version mismatch" Error Detail
# messages are cleared now that everything is in sync, and the zh-cn
# handback datetime is refreshed. Column widths are widened/narrowed to
# fit the new text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (same text as the
# per-locale "Status" column below)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
# Latest Handback DateTime refreshed
$zhcn.Range("K2").Value = "2016-09-07 07:58:47"
$zhcn.Range("K3").Value = "2016-09-07 07:58:47"
# Error Detail no longer applicable - handback is in sync now
$zhcn.Range("P2").Value = ""
$zhcn.Range("P3").Value = ""

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus
# Latest Handoff File / Latest Handback File now the real de-de xliff
$dedeFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("G2").Value = $dedeFile
$dede.Range("G3").Value = $dedeFile
$dede.Range("J2").Value = $dedeFile
$dede.Range("J3").Value = $dedeFile
# Latest Handback DateTime refreshed
$dede.Range("K2").Value = "2016-09-07 07:59:07"
$dede.Range("K3").Value = "2016-09-07 07:59:07"
# Error Detail no longer applicable - handback is in sync now
$dede.Range("P2").Value = ""
$dede.Range("P3").Value = ""

# ---------------------------------------------------------------------
# Column widths - widen the Status columns to fit the longer text,
# narrow the Error Detail column now that it's empty.
# ---------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333334

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333334
